$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each Price/Volume cell is stored as text in the source data (e.g. "59.917.62",
# "0.998", "  -4.06%  "). Excel auto-converts plain numeric-looking strings typed
# into Value to actual numbers, so we prefix with a literal single-quote (forces
# text entry / quote-prefix) and then reset the cell Style back to "Normal" so no
# stray number-format/quote-prefix style is left behind - matching the original
# plain (unstyled) text cells.

$ws.Range("D2").Value = "'59.917.62"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -4.06%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.492.68"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -5.16%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.21%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'543.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.54%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'147.05"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -4.84%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.996"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.36%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.581"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -1.09%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'2.521.27"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -4.17%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -3.41%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -0.88%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'5.43"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.28%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  -1.82%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'2.936.77"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -5.20%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'24.64"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -4.42%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'59.685.94"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -4.30%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  -2.73%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.520.96"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -4.21%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'11.46"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -1.54%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'4.37"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -3.55%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'327.39"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -3.81%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.996"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.42%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'5.81"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -4.78%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'61.46"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -2.41%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.449"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -10.10%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +0.79%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -3.36%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'7.85"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -2.47%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'1.30"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -2.03%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.0₃0793"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -5.13%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'6.88"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -3.31%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'1.83"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -3.69%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.998"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -0.05%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'158.78"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -1.55%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.45"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +1.23%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'18.93"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -1.57%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'4.53"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -4.36%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'1.74"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.14%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'6.06"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -1.28%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'314.30"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -6.89%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'36.79"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -2.97%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'3.78"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -3.47%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.837"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -7.88%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.996"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.27%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.607"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.68%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'10.79"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -1.77%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'126.84"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.50%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.0533"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -2.50%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.0944"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -1.87%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0232"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -2.48%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'18.70"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -5.46%  "
$ws.Range("E51").Style = "Normal"
